$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.498.77'
$ws.Range("E2").Value = '  -1.63%  '

$ws.Range("D3").Value = '1.650.49'
$ws.Range("E3").Value = '  -3.45%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.69%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9994'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.18%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3650'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.68%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '46.48'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.80%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3245'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.90%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.122'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.87%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07016'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -7.00%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9997'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.13%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.953'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.74%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.33'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -9.23%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.595'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.84%  '

$ws.Range("D16").Value = '1.652.97'
$ws.Range("E16").Value = '  -3.13%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001040'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -8.15%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06587'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.01%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9984'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.06%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '78.50'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -7.01%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.921'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -7.41%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.65'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -9.49%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.50'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.72%  '

$ws.Range("D24").Value = '24.481.71'
$ws.Range("E24").Value = '  -1.76%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.468'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.10%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.315'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -17.28%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '146.71'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.20%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.50'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -9.52%  '

$ws.Range("D29").Value = '1.833.55'
$ws.Range("E29").Value = '  -3.21%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.13'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -6.84%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.183'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.18%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.069'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.89%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.689'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -17.18%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08425'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.25%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.662'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.46%  '

$ws.Range("E36").Value = '  -12.19%  '

$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.270'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.73%  '

$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.178'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -8.31%  '

$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06005'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -9.86%  '

$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02216'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -8.17%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2066'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -7.48%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.063'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -12.27%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9987'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.09%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5875'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -9.10%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.789'
$ws.Range("D45").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.64'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -8.78%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5601'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -9.08%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '123.38'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.03%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.940'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -8.87%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06916'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.55%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.187'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.86%  '
